$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column AV (48) needs the same column width as the rest of the data columns (12 chars)
$ws.Columns.Item(48).ColumnWidth = 11.17

# --- Header cell AV1 ---
# Copy AU1 formatting (style index 1: Meiryo font, no fill, General number format)
$ws.Cells.Item(1, 47).Copy($ws.Cells.Item(1, 48))

# Write the literal text "2024/10/26" without Excel re-interpreting it as a date:
# stage it in a scratch cell that is explicitly Text-formatted, paste the VALUE only
# into AV1 (so AV1 keeps the style copied above), then wipe the scratch cell.
$scratch = $ws.Cells.Item(200, 200)
$scratch.NumberFormat = "@"
$scratch.Value = "2024/10/26"
$scratch.Copy()
$ws.Cells.Item(1, 48).PasteSpecial(-4163)
$scratch.Clear()

# --- Data cells AV2:AV53 ---
# The sheet colour-codes values via baked-in cell styles (no real conditional formatting):
#   value <  125          -> style 2 (yellow fill)  : reference cell D2
#   125 <= value < 140    -> style 3 (blue fill)    : reference cell N2
#   value >= 140          -> style 1 (no fill)      : reference cell A2
$styleRefs = @{ 1 = $ws.Cells.Item(2, 1); 2 = $ws.Cells.Item(2, 4); 3 = $ws.Cells.Item(2, 14) }

$avData = @(
    @{ Row = 2; Style = 3; Value = 135.2 },
    @{ Row = 3; Style = 1; Value = 183.5 },
    @{ Row = 4; Style = 1; Value = 149.6 },
    @{ Row = 5; Style = 1; Value = 146.7 },
    @{ Row = 6; Style = 1; Value = 195.3 },
    @{ Row = 7; Style = 1; Value = 182.3 },
    @{ Row = 8; Style = 3; Value = 135.8 },
    @{ Row = 9; Style = 1; Value = 283.6 },
    @{ Row = 10; Style = 1; Value = 165.1 },
    @{ Row = 11; Style = 3; Value = 139 },
    @{ Row = 12; Style = 1; Value = 268 },
    @{ Row = 13; Style = 1; Value = 202.1 },
    @{ Row = 14; Style = 1; Value = 151.9 },
    @{ Row = 15; Style = 1; Value = 146.1 },
    @{ Row = 16; Style = 1; Value = 144.2 },
    @{ Row = 17; Style = 1; Value = 188.5 },
    @{ Row = 18; Style = 1; Value = 162.4 },
    @{ Row = 19; Style = 1; Value = 143.2 },
    @{ Row = 20; Style = 1; Value = 165.3 },
    @{ Row = 21; Style = 1; Value = 159.6 },
    @{ Row = 22; Style = 1; Value = 181.5 },
    @{ Row = 23; Style = 2; Value = 116.9 },
    @{ Row = 24; Style = 1; Value = 205.3 },
    @{ Row = 25; Style = 1; Value = 166.7 },
    @{ Row = 26; Style = 1; Value = 171.7 },
    @{ Row = 27; Style = 3; Value = 132.9 },
    @{ Row = 28; Style = 3; Value = 137.3 },
    @{ Row = 29; Style = 3; Value = 134.2 },
    @{ Row = 30; Style = 1; Value = 155.4 },
    @{ Row = 31; Style = 3; Value = 131 },
    @{ Row = 32; Style = 3; Value = 131.3 },
    @{ Row = 33; Style = 1; Value = 252.6 },
    @{ Row = 34; Style = 1; Value = 140 },
    @{ Row = 35; Style = 1; Value = 155.2 },
    @{ Row = 36; Style = 1; Value = 154.8 },
    @{ Row = 37; Style = 3; Value = 135.8 },
    @{ Row = 38; Style = 1; Value = 151.2 },
    @{ Row = 39; Style = 1; Value = 161.5 },
    @{ Row = 40; Style = 1; Value = 182.7 },
    @{ Row = 41; Style = 1; Value = 141.1 },
    @{ Row = 42; Style = 1; Value = 157.6 },
    @{ Row = 43; Style = 1; Value = 200.1 },
    @{ Row = 44; Style = 2; Value = 120.1 },
    @{ Row = 45; Style = 1; Value = 146.2 },
    @{ Row = 46; Style = 1; Value = 171.7 },
    @{ Row = 47; Style = 1; Value = 342.4 },
    @{ Row = 48; Style = 1; Value = 188.1 },
    @{ Row = 49; Style = 3; Value = 136.4 },
    @{ Row = 50; Style = 1; Value = 147.3 },
    @{ Row = 51; Style = 1; Value = 162.1 },
    @{ Row = 52; Style = 1; Value = 156.1 },
    @{ Row = 53; Style = 1; Value = 172.1 }
)

foreach ($item in $avData) {
    $styleRefs[$item.Style].Copy($ws.Cells.Item($item.Row, 48))
    $ws.Cells.Item($item.Row, 48).Value = $item.Value
}

Write-Output "done"